$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) " do not bear fruit so beautiful." -> " do not bear fruit so beautifully."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("do not bear fruit so beautiful.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "do not bear fruit so beautifully.", 2)

# ---------------------------------------------------------------------------
# 2) "They are very difficult to " + "bore"  ->  "...to <tl>" + "bore</tl>"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("They are very difficult to ")
$rng.Collapse(0)
$rng.InsertAfter("<tl>")

$rng = $d.Content
$rng.Find.Execute(" when they are of ")
$rng.Collapse(1)
$rng.InsertBefore("</tl>")

# ---------------------------------------------------------------------------
# 3a) "one ought" -> "one <del><fr>z</fr></del> ought"  (insertion inside the run)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("one ought")
$p = $rng.Start + 4
$ins = $d.Range($p, $p)
$ins.InsertBefore("<del><fr>z</fr></del> ")

# ---------------------------------------------------------------------------
# 3b) " not to push the " + "borer"  ->  "... the <tl>" + "borer</tl>"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" not to push the ")
$rng.Collapse(0)
$rng.InsertAfter("<tl>")

$rng = $d.Content
$rng.Find.Execute(" with too great a force")
$rng.Collapse(1)
$rng.InsertBefore("</tl>")

# ---------------------------------------------------------------------------
# 4) "</del> " + "muzzle of the cannon "  ->  "</del> <add>" + "muzzle of the cannon</add> "
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("</del> ")
$rng.Collapse(0)
$rng.InsertAfter("<add>")

$rng = $d.Content
$rng.Find.Execute("muzzle of the cannon")
$rng.Collapse(0)
$rng.InsertAfter("</add>")

# ---------------------------------------------------------------------------
# 5) " in order to avoid the necessity of " + "boring"  ->  "...of <tl>" + "boring</tl>"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" in order to avoid the necessity of ")
$rng.Collapse(0)
$rng.InsertAfter("<tl>")

$rng = $d.Content
$rng.Find.Execute("boring")
$rng.Collapse(0)
$rng.InsertAfter("</tl>")

# ---------------------------------------------------------------------------
# 6) ". The " + "borer"  ->  ". The <tl>" + "borer</tl>"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(". The ")
$rng.Collapse(0)
$rng.InsertAfter("<tl>")

$rng = $d.Content
$rng.Find.Execute(" should be neither like a")
$rng.Collapse(1)
$rng.InsertBefore("</tl>")

# ---------------------------------------------------------------------------
# 7) " round like a " + "wimble"  ->  "...like a <tl>" + "wimble</tl>"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" round like a ")
$rng.Collapse(0)
$rng.InsertAfter("<tl>")

$rng = $d.Content
$rng.Find.Execute(". If it breaks, ")
$rng.Collapse(1)
$rng.InsertBefore("</tl>")

# ---------------------------------------------------------------------------
# 8) "left-top"  ->  "left-" (same run/style) + "middle" (new run, plain style)
#    Copy an existing plain-style run ("bore") to get Word to mint a run with
#    the bare <w:rtl w:val="0"/> rPr, then retarget its text to "middle".
# ---------------------------------------------------------------------------
$rngSrc = $d.Content
$rngSrc.Find.Execute("bore")
$rngSrc.Copy()

$rng = $d.Content
$rng.Find.Execute("left-top")
$rng.Text = "left-"
$insPoint = $rng.End
$rngIns = $d.Range($insPoint, $insPoint)
$rngIns.Paste()

$rngMid = $d.Range($insPoint, $insPoint + 4)
$rngMid.Text = "middle"

# ---------------------------------------------------------------------------
# 9) " then turn the place of the touch-hole"
#    -> " then <del><fr>l</fr></del> turn the place of the touch-hole"
#    (insertion inside the run, no run split)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" then turn the place of the touch-hole")
$p = $rng.Start + 6
$ins = $d.Range($p, $p)
$ins.InsertBefore("<del><fr>l</fr></del> ")

# ---------------------------------------------------------------------------
# 10) "One needs to bore gently &" -> "One needs to <tl>bore</tl> gently &"
#     (insertion inside the run, no run split)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("One needs to bore gently")
$p1 = $rng.Start + 13
$ins1 = $d.Range($p1, $p1)
$ins1.InsertBefore("<tl>")

$p2 = $rng.Start + 13 + 4 + 4
$ins2 = $d.Range($p2, $p2)
$ins2.InsertBefore("</tl>")

# ---------------------------------------------------------------------------
# 11) " piece is of metal." -> " piece is of <m>metal</m>."
#     (insertion inside the run, no run split)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" piece is of metal.")
$q1 = $rng.Start + 13
$ins3 = $d.Range($q1, $q1)
$ins3.InsertBefore("<m>")

$q2 = $rng.Start + 13 + 3 + 5
$ins4 = $d.Range($q2, $q2)
$ins4.InsertBefore("</m>")
